$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.473.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.024.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.65%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.439"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("E11").Value = "  +3.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.545.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("E13").Value = "  +1.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.20%  "

$ws.Range("E15").Value = "  +7.52%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.512.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.45%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.031.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("E19").Value = "  +2.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("E23").Value = "  +4.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0929"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.19%  "

$ws.Range("E30").Value = "  +2.45%  "

$ws.Range("E31").Value = "  -3.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.66%  "

$ws.Range("E35").Value = "  +4.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "24.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.18%  "

$ws.Range("E38").Value = "  +1.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.060.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.87%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.303.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.655"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.01%  "

$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.991"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.19%  "

$ws.Range("E48").Value = "  +1.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.78%  "

$ws.Range("E50").Value = "  -7.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0893"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.95%  "
